# Update "paises.xlsx" (Pais worksheet) with the latest COVID-19 country stats
# and refresh the "last updated" timestamp, per the Aug 4 2020 17:58 data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: refresh "last updated" timestamp ---
$ws.Range("A1").Value = 'Datos actualizados a 4 de Agosto de 2020 a las 17:58'

# --- Countries that changed rank order (new leader gets fresh stats, the
#     country it overtakes keeps its own previous stats one row down) ---
# Ucrania / Republica Dominicana swap rank positions 41-42
$ws.Range("A37").Value = 'Republica Dominicana'
$ws.Range("B37").Value = 74295
$ws.Range("C37").Value = 1178
$ws.Range("D37").Value = 38824
$ws.Range("E37").Value = 34258
$ws.Range("G37").Value = 30
$ws.Range("H37").Value = 1213
$ws.Range("A38").Value = 'Ucrania'
$ws.Range("B38").Value = 74219
$ws.Range("C38").Value = 1061
$ws.Range("D38").Value = 40613
$ws.Range("E38").Value = 31842
$ws.Range("G38").Value = 26
$ws.Range("H38").Value = 1764

# Etiopia / Costa Rica swap rank positions 74-75
$ws.Range("A70").Value = 'Etiopia'
$ws.Range("B70").Value = 19877
$ws.Range("C70").Value = 588
$ws.Range("D70").Value = 8240
$ws.Range("E70").Value = 11294
$ws.Range("G70").Value = 7
$ws.Range("H70").Value = 343
$ws.Range("A71").Value = 'Costa Rica'
$ws.Range("B71").Value = 19402
$ws.Range("D71").Value = 4689
$ws.Range("E71").Value = 14542
$ws.Range("H71").Value = 171

# Grecia / Guinea Ecuatorial swap rank positions 107-108
$ws.Range("A103").Value = 'Grecia'
$ws.Range("B103").Value = 4855
$ws.Range("C103").Value = 118
$ws.Range("D103").Value = 1374
$ws.Range("E103").Value = 3272
$ws.Range("H103").Value = 209
$ws.Range("A104").Value = 'Guinea Ecuatorial'
$ws.Range("B104").Value = 4821
$ws.Range("D104").Value = 2182
$ws.Range("E104").Value = 2556
$ws.Range("H104").Value = 83

# Trinidad yTobago / Gibraltar swap rank positions 183-184
$ws.Range("A179").Value = 'Trinidad yTobago'
$ws.Range("B179").Value = 192
$ws.Range("C179").Value = 10
$ws.Range("D179").Value = 135
$ws.Range("E179").Value = 49
$ws.Range("H179").Value = 8
$ws.Range("A180").Value = 'Gibraltar'
$ws.Range("B180").Value = 189
$ws.Range("C180").Value = 1
$ws.Range("D180").Value = 183
$ws.Range("E180").Value = 6
$ws.Range("H180").Value = 0

# --- Remaining countries: daily totals refreshed in place ---
$ws.Range("B4").Value = 4873103
$ws.Range("C4").Value = 10929
$ws.Range("D4").Value = 2449120
$ws.Range("E4").Value = 2264614
$ws.Range("G4").Value = 441
$ws.Range("H4").Value = 159369
$ws.Range("B5").Value = 2755081
$ws.Range("C5").Value = 3416
$ws.Range("E5").Value = 747981
$ws.Range("G5").Value = 79
$ws.Range("H5").Value = 94781
$ws.Range("B6").Value = 1901334
$ws.Range("C6").Value = 46003
$ws.Range("D6").Value = 1278084
$ws.Range("E6").Value = 583463
$ws.Range("G6").Value = 816
$ws.Range("H6").Value = 39787
$ws.Range("B15").Value = 306293
$ws.Range("C15").Value = 670
$ws.Range("B18").Value = 248419
$ws.Range("C18").Value = 190
$ws.Range("D18").Value = 200766
$ws.Range("E18").Value = 12482
$ws.Range("G18").Value = 5
$ws.Range("H18").Value = 35171
$ws.Range("D23").Value = 82166
$ws.Range("E23").Value = 78835
$ws.Range("B25").Value = 117210
$ws.Range("C25").Value = 179
$ws.Range("D25").Value = 101839
$ws.Range("E25").Value = 6420
$ws.Range("G25").Value = 4
$ws.Range("H25").Value = 8951
$ws.Range("B65").Value = 25814
$ws.Range("C65").Value = 332
$ws.Range("E65").Value = 7062
$ws.Range("G65").Value = 10
$ws.Range("H65").Value = 810
$ws.Range("B75").Value = 17114
$ws.Range("C75").Value = 106
$ws.Range("D75").Value = 11808
$ws.Range("E75").Value = 4920
$ws.Range("B119").Value = 2834
$ws.Range("C119").Value = 6
$ws.Range("E119").Value = 299
$ws.Range("B177").Value = 227
$ws.Range("C177").Value = 2
$ws.Range("E177").Value = 35
